# PAS-6576: fix typos in VIN upload test excel file
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: BI_SYMBOL/PD_SYMBOL/UM_SYMBOL/MP_SYMBOL changed from "X" to "A"
$ws.Range("AE3:AH3").Value = "A"

# Row 4: BI_SYMBOL/PD_SYMBOL/UM_SYMBOL/MP_SYMBOL changed from "I" to "X"
$ws.Range("AE4:AH4").Value = "X"

# Update the selected cell/range shown in the sheet view
$ws.Range("H18").Select()
